$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(224,122,209,107,169,251,273,71,132)
    ,@(64,26,76,29,52,62,58,20,36)
    ,@(11,4,7,2,11,12,11,4,12)
    ,@(6,1,8,1,4,6,2,1,8)
    ,@(11,6,12,4,7,5,9,0,2)
    ,@(1,2,8,2,2,2,6,3,5)
    ,@(19,7,18,3,12,13,13,4,1)
    ,@(2,0,3,0,1,3,2,2,2)
    ,@(6,1,10,1,0,10,3,0,0)
    ,@(0,0,0,0,0,0,0,0,0)
    ,@(8,5,10,16,15,11,12,6,6)
    ,@(62,50,56,20,52,31,89,33,40)
    ,@(16,15,1,1,17,8,11,10,2)
    ,@(10,6,11,2,3,6,8,9,3)
    ,@(5,16,9,0,12,4,21,6,9)
    ,@(6,6,10,7,6,6,13,8,8)
    ,@(2,1,9,3,5,5,9,0,12)
    ,@(8,6,5,7,9,2,10,0,3)
    ,@(15,0,11,0,0,0,17,0,3)
    ,@(32,12,24,22,26,32,44,4,12)
    ,@(2,4,6,3,3,2,9,3,5)
    ,@(1,0,2,0,2,3,1,1,1)
    ,@(11,4,6,7,8,6,10,0,2)
    ,@(7,4,6,6,10,9,11,0,3)
    ,@(11,0,4,6,3,12,13,0,1)
    ,@(25,7,28,17,15,64,32,7,17)
    ,@(3,0,6,0,0,5,5,0,2)
    ,@(2,0,3,0,0,17,0,1,1)
    ,@(2,1,4,1,2,14,16,2,4)
    ,@(10,2,7,9,10,15,9,4,5)
    ,@(8,4,8,7,3,13,2,0,5)
    ,@(27,15,18,15,16,31,29,4,13)
    ,@(5,2,4,5,7,6,8,1,1)
    ,@(8,7,3,5,4,2,3,0,5)
    ,@(0,0,0,0,0,0,0,0,0)
    ,@(7,3,8,3,2,9,10,0,5)
    ,@(3,2,1,0,2,4,4,1,0)
    ,@(4,1,2,2,1,10,4,2,2)
    ,@(0,0,0,0,0,0,0,0,0)
    ,@(14,12,7,4,8,31,21,3,14)
    ,@(6,6,0,0,0,10,2,1,3)
    ,@(0,0,0,0,0,0,0,0,0)
    ,@(1,2,4,2,3,2,8,1,1)
    ,@(1,1,2,1,3,15,7,0,10)
    ,@(6,3,1,1,2,4,4,1,0)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $row[$c]
    }
}
